$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# D32 previously held a hard-coded value; replace it with a formula that
# references the prior row's running total (F31), matching the pattern
# already used by D3:D29 ( =F<prev row> ).
$ws.Range("D32").Formula = "=F31"

# D33 likewise should reference F32.
$ws.Range("D33").Formula = "=F32"

$wb.Application.Calculate()
